$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: new product id, new product name, new price
$ws.Range("A2").Value = 5151
$ws.Range("B2").Value = "ريد بل - 250 مل"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1065
$ws.Range("E2").Value = "YES"

# Add new row 3 with the previous product under a new name
$ws.Range("A3").Value = 7630
$ws.Range("B3").Value = "فيورى جولد - 400 مل"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 205
$ws.Range("E3").Value = "YES"
